# The deck's date-footer placeholder ("today's date" field) was cached as
# 3/28/22 on the Slide Master and on every Slide Layout; the edit simply
# refreshes that cached display text to 3/31/22 everywhere it appears
# (ppPlaceholderDate, PlaceholderFormat.Type = 16), exactly as PowerPoint
# itself re-stamps the datetimeFigureOut field's <a:t> on open/save.

$p = $ppt.ActivePresentation
$newDate = "3/31/22"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
        }

        if (-not $isDatePlaceholder) {
            if ($sh.Name -like "Date Placeholder*") {
                $isDatePlaceholder = $true
            }
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout's date placeholder.
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Belt-and-suspenders: also sweep any slide that carries its own
# (overridden) date placeholder instead of inheriting from its layout.
for ($S = 1; $S -le $p.Slides.Count; $S++) {
    Update-DatePlaceholder $p.Slides.Item($S).Shapes
}

Write-Output "Updated date placeholders to $newDate"
